# Generate Report for Archive
# - Update the localization status text from "Ready for handoff" to
#   "In Translation" on every sheet/cell that currently shows it
#   (Overview!E2/F2 summary columns, and the per-language Status column
#   C2 on the "zh-cn" and "de-de" sheets).
# - Narrow the now-shorter "Status" columns to match (Overview cols E/F,
#   and column C on the zh-cn / de-de sheets).

$wb = $excel.ActiveWorkbook

# was "Ready for handoff"
$newStatus = "In Translation"

# Target OOXML column "width" the commit lands on is 13.4101845877511
# characters. This engine's ColumnWidth setter quantises to 1/6-character
# (pixel) steps before re-deriving the stored width, so feed it the input
# that lands on the closest achievable step (13.333333333333334).
$newColumnWidth = 12.5

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
